$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new SVR kernel parameter headers
$ws.Range("K1").Value = "svr_kernel_scale"
$ws.Range("L1").Value = "svr_epsilon"
$ws.Range("M1").Value = "svr_box_constraint"

# Add the corresponding values
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 5

# Update selection to match target state
$ws.Range("L9").Select()
